$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.94"
$ws.Range("E2").Value = "'1.88%"
$ws.Range("D3").Value = "'27.19"
$ws.Range("E3").Value = "'1.91%"
$ws.Range("D4").Value = "'4.728"
$ws.Range("E4").Value = "'5.05%"
$ws.Range("D5").Value = "'0.06080"
$ws.Range("E5").Value = "'3.46%"
$ws.Range("D6").Value = "'6.667"
$ws.Range("E6").Value = "'0.89%"
$ws.Range("D7").Value = "'0.8477"
$ws.Range("E7").Value = "'-0.36%"
$ws.Range("D8").Value = "'0.9217"
$ws.Range("E8").Value = "'-0.70%"
$ws.Range("E9").Value = "'2.14%"
$ws.Range("D10").Value = "'0.04941"
$ws.Range("E10").Value = "'8.74%"
$ws.Range("D11").Value = "'0.07092"
$ws.Range("E11").Value = "'0.77%"
$ws.Range("D12").Value = "'0.03135"
$ws.Range("E12").Value = "'2.22%"
$ws.Range("D13").Value = "'0.09077"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("D14").Value = "'0.001537"
$ws.Range("E14").Value = "'0.39%"
$ws.Range("D15").Value = "'0.0006111"
$ws.Range("E15").Value = "'0.83%"
$ws.Range("D16").Value = "'0.006160"
$ws.Range("E16").Value = "'-0.60%"
$ws.Range("D17").Value = "'3.450"
$ws.Range("E17").Value = "'-0.92%"
$ws.Range("D18").Value = "'3.155"
$ws.Range("E18").Value = "'-0.59%"
$ws.Range("E19").Value = "'-2.22%"
$ws.Range("E20").Value = "'2.59%"
$ws.Range("E21").Value = "'0.87%"
$ws.Range("D22").Value = "'4.099"
$ws.Range("E22").Value = "'4.39%"
$ws.Range("D23").Value = "'0.04248"
$ws.Range("E23").Value = "'-0.16%"
$ws.Range("E24").Value = "'0.11%"
$ws.Range("E25").Value = "'-8.91%"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("D27").Value = "'0.0001575"
$ws.Range("E27").Value = "'3.37%"
$ws.Range("D40").Value = "'0.03877"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'1.35%"
$ws.Range("D42").Value = "'0.004127"
$ws.Range("E42").Value = "'5.86%"
$ws.Range("D43").Value = "'0.01634"
$ws.Range("E43").Value = "'17.88%"
$ws.Range("E44").Value = "'-9.12%"
$ws.Range("D45").Value = "'0.00005334"
$ws.Range("E45").Value = "'-0.32%"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("E47").Value = "'1.25%"
$ws.Range("D48").Value = "'0.1321"
$ws.Range("E48").Value = "'-47.56%"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E50").Value = "'0.01%"
